$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "department" column (E), filled for existing rows first ---
$ws.Cells.Item(1,5).Value = "department"
$ws.Cells.Item(2,5).Value = "Econ"
$ws.Cells.Item(3,5).Value = "Econ"

# --- New row 4: Jane Doe candidate record ---
$ws.Cells.Item(4,1).Value = "Jane Doe"
$ws.Cells.Item(4,2).Value = "This is a really long JMP title to test character limits: Evidence from a randomized control trial"

$ws.Cells.Item(4,3).Value = "https://www.twitter.com/"
$ws.Hyperlinks.Add($ws.Cells.Item(4,3), "https://www.twitter.com/")
$ws.Cells.Item(4,3).Style = "Hyperlink"

$ws.Cells.Item(4,4).Value = "University of Phoenix"
$ws.Cells.Item(4,5).Value = "Econ"

# --- New "posted" / "qflag" quality-control columns (F, G) ---
$ws.Cells.Item(1,6).Value = "posted"
$ws.Cells.Item(1,7).Value = "qflag"

$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = 0
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 0

# --- Column E width to match the other data columns ---
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# --- Move the active selection to G4, matching the refreshed QC layout ---
[void]$ws.Range("G4").Select()
